$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 129 - this shifts rows 129:165 down to 130:166
# and keeps all their existing data/formatting intact.
$ws.Rows.Item(129).Insert()

# Populate the newly inserted row 129 with the new record's data. All the
# "constant" columns (A,B,C,E,F,G,H,I,J,K,Q,T) carry the same values as
# every other row in this sheet.
$ws.Range("A129").Value = 9
$ws.Range("B129").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C129").Value = "Metropolitana"
$ws.Range("D129").Value = 44588
$ws.Range("E129").Value = 13
$ws.Range("F129").Value = "Fruta"
$ws.Range("G129").Value = 100101
$ws.Range("H129").Value = "Berries"
$ws.Range("I129").Value = 100101001
$ws.Range("J129").Value = "Arándano (blue)"
$ws.Range("K129").Value = "Sin especificar"
$ws.Range("L129").Value = "Primera"
$ws.Range("M129").Value = 350
$ws.Range("N129").Value = 4000
$ws.Range("O129").Value = 4000
$ws.Range("P129").Value = 4000
$ws.Range("Q129").Value = "$/bandeja 2 kilos"
$ws.Range("R129").Value = "Región de O'Higgins"
$ws.Range("S129").Value = 2000
$ws.Range("T129").Value = 2
